# Update "想去人数" (F column) counts for several events on both the
# "展览" sheet and the "全部类型" sheet, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 875
$ws1.Range("F4").Value = 18
$ws1.Range("F8").Value = 4615
$ws1.Range("F10").Value = 38
$ws1.Range("F11").Value = 506
$ws1.Range("F12").Value = 464
$ws1.Range("F13").Value = 19
$ws1.Range("F16").Value = 2650
$ws1.Range("F17").Value = 382
$ws1.Range("F18").Value = 80
$ws1.Range("F20").Value = 60
$ws1.Range("F21").Value = 2303
$ws1.Range("F22").Value = 93
$ws1.Range("F24").Value = 30
$ws1.Range("F25").Value = 162
$ws1.Range("F28").Value = 227
$ws1.Range("F29").Value = 37

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 875
$ws4.Range("F4").Value = 18
$ws4.Range("F9").Value = 4615
$ws4.Range("F11").Value = 38
$ws4.Range("F12").Value = 506
$ws4.Range("F13").Value = 464
$ws4.Range("F14").Value = 19
$ws4.Range("F17").Value = 2650
$ws4.Range("F18").Value = 382
$ws4.Range("F19").Value = 80
$ws4.Range("F21").Value = 60
$ws4.Range("F22").Value = 2303
$ws4.Range("F23").Value = 93
$ws4.Range("F25").Value = 30
$ws4.Range("F26").Value = 162
$ws4.Range("F29").Value = 227
$ws4.Range("F30").Value = 37
